$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "TimeLine and Cost" -> "Timeline and Cost" (fix capitalisation)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("TimeLine and Cost", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Timeline and Cost", 2)

# ---------------------------------------------------------------------
# Change 2: rewrite the MVP timeline/cost paragraph with the new
# estimate and longer explanation, switching it to the Calibri/222222
# body-text formatting used elsewhere in the document.
# ---------------------------------------------------------------------
$old2 = "The MVP should be delivered in 2 weeks with a cost of around 50K. "
$new2 = "The MVP should be delivered in 2 weeks with a cost of around `$70K-`$75K. This is an approximation and details will be provided after subsequent planning sessions. Based on customer feedback more features can be added, and additional cost and timeline will be projected. "

$styleSrc = $d.Content
$styleSrc.Find.Execute("Employees are a company", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

$target = $d.Content
$target.Find.Execute($old2, $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

$target.FormattedText = $styleSrc.FormattedText
$retext = $d.Range($target.Start, $target.Start + $styleSrc.Text.Length)
$retext.Text = $new2

# ---------------------------------------------------------------------
# Change 3: merge the runs split by the old page-break position so the
# phrase reads as one run again ("...MVP works well with the...").
# ---------------------------------------------------------------------
$merge = $d.Content
$merge.Find.Execute(" If the MVP works well with the ", $false, $false, $false, $false, $false, `
    $true, 1, $false, " If the MVP works well with the ", 1)

# ---------------------------------------------------------------------
# Change 3 (cont.): the page break now falls at the top of the next
# paragraph, so mark it as the last rendered page break there.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Employees are a company", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$insertionPoint = $d.Range($anchor.Start, $anchor.Start)
$breakXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="222222"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/></w:r></w:p>'
$insertionPoint.InsertXML($breakXml)

Write-Output "done"
